$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column A: existing Title/Path columns shift from A,B to B,C ---
$ws.Columns("A:A").Insert()

# --- Header row (row 1) ---
$ws.Range("A1").Value = "Source"
$ws.Range("E1").Value = "RelativePath"
# B1 / C1 already hold "Title" / "StratScreenshot" after the column shift.

# --- Column A: "Core" source tag for every data row ---
$ws.Range("A2:A12").Value = "Core"

# --- Column E: shared relative-path prefix for every data row ---
$ws.Range("E2:E12").Value = "strats-sources\assets\v10_screenshots"

# --- Column B: strat titles for the newly added rows (6-12) ---
$ws.Range("B6").Value = "Defi Epique"
$ws.Range("B7").Value = "Attaque de Char"
$ws.Range("B8").Value = "Tir en Etat d'alerte"
$ws.Range("B9").Value = "Arrivee precipitee"
$ws.Range("B10").Value = "A Terre"
$ws.Range("B11").Value = "Ecran de Fumee"
$ws.Range("B12").Value = "Intervention Heroique"

# --- Column D: screenshot filenames for every data row ---
$ws.Range("D2").Value = "00_Core_RelanceCommandement.png"
$ws.Range("D3").Value = "01_Core_ContreOffensive.png"
$ws.Range("D4").Value = "02_Core_DefiEpique.png"
$ws.Range("D5").Value = "03_Core_CourageInsense.png"
$ws.Range("D6").Value = "04_Core_Grenade.png"
$ws.Range("D7").Value = "05_Core_AttaqueDeChar.png"
$ws.Range("D8").Value = "06_Core_TirEtatAlerte.png"
$ws.Range("D9").Value = "07_Core_ArriveePrecipitee.png"
$ws.Range("D10").Value = "08_Core_ATerre.png"
$ws.Range("D11").Value = "09_Core_EcranFumee.png"
$ws.Range("D12").Value = "10_Core_InterventionHeroique.png"

# --- Column C: formula combining E (path) and D (filename) ---
$ws.Range("C2").Formula = '=(E2 & "\" & D2)'
$ws.Range("C3:C12").Formula = '=(E3 & "\" & D3)'

# --- Column widths for the two new columns ---
$ws.Columns("D:D").ColumnWidth = 71
$ws.Columns("E:E").ColumnWidth = 35

# --- View: zoom + final selection, matching the saved author state ---
$ws.Application.ActiveWindow.Zoom = 85
$null = $ws.Range("C12").Select()
